# Add a new "Player Info" worksheet as the first sheet in the workbook and
# update the two existing sheets so their MATCH_CARD_LINK column becomes a
# MATCH_CODE column holding just the bare numeric match code instead of the
# full scorecard URL.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Player Info" sheet in front of everything else ----
# NOTE: worksheet references are position-based, so grab/rename the new
# sheet *before* looking up the other sheets by name - re-resolving "before"
# and "bowling" only after the insert keeps their handles pointed at the
# right tab.
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

# Header row - bold, centered, thin-bordered (matches the other sheets)
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row. The ID looks numeric but must be stored as text (like every
# other "numeric-looking" column in this workbook), so force text with a
# leading quote and then drop back to the normal (unstyled) cell style.
$playerInfo.Range("A2").Value = "'6859"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Lokumarakkalage Dilshan Madushanka"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Left Arm Medium Fast"

# --- 2. Update the "ODI Batting" sheet -------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").Value = "'4687"
$batting.Range("D2").Style = "Normal"

# --- 3. Update the "ODI Bowling" sheet -------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "'4687"
$bowling.Range("B2").Style = "Normal"
